# Update Name of Algo
# Apply updated KNN-imputed values to the result_data_KNN worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = 8.129000000000001
    "D4"  = -7.968999999999999
    "B6"  = 5.689
    "B7"  = 5.732
    "D9"  = -7.898000000000001
    "D12" = -7.072
    "B16" = 4.947
    "D17" = -8.191999999999998
    "D18" = -8.365
    "D19" = -7.967000000000001
    "B20" = 8.132000000000001
    "D20" = -7.897
    "D26" = -7.351000000000001
    "B28" = 5.497999999999999
    "B29" = 5.281
    "D31" = -7.817000000000002
    "B32" = 6.48
    "D39" = -7.632
    "B40" = 9.554999999999998
    "D40" = -8.019
    "D41" = -7.997999999999999
    "D42" = -8.036
    "D43" = -7.777000000000001
    "B46" = 5.613
    "D47" = -7.639
    "D48" = -7.651999999999999
    "B51" = 5.548
    "B52" = 5.598
    "B57" = 5.040999999999999
    "B59" = 4.988
    "B62" = 5.311999999999999
    "D63" = -6.923
    "D64" = -7.263
    "B66" = 5.211
    "B73" = 7.145
    "B74" = 9.164999999999999
    "D76" = -7.748
    "D81" = -8.013
    "D89" = -8.291999999999998
    "B92" = 5.178
    "D94" = -7.640000000000001
    "B100" = 6.073
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
